# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gets a new blank column inserted before
# column N (pushing the old "Late" / "Outstanding" columns one slot to the
# right: N->O, O->P, P->Q), and the sheet becomes the active sheet/tab with
# the selection moved to R5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N, shifting existing N:P -> O:Q.
$ws.Columns("N:N").Insert()

# New column gets a plain (non bestFit) width of 10.
$ws.Columns("N:N").ColumnWidth = 9.140625

# Make "Repayment Schedule" the active sheet/tab, with R5 selected.
$ws.Activate()
$null = $ws.Range("R5").Select()
